$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 5383.65
$ws.Range("J40").Value2 = 6721.778
$ws.Range("L40").Value2 = 6721.778
$ws.Range("N40").Value2 = -7071.778
$ws.Range("H51").Value2 = 4374.5835
$ws.Range("I51").Value2 = 3271.2856
$ws.Range("J51").Value2 = 4828.8823
$ws.Range("K51").Value2 = 3271.2856
$ws.Range("L51").Value2 = 4828.8823
$ws.Range("M51").Value2 = -2787.2856
$ws.Range("N51").Value2 = -5796.8823
$ws.Range("H92").Value2 = 2241.4167
$ws.Range("I92").Value2 = 2241.4167
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 2241.4167
$ws.Range("L92").Value2 = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value2 = -993.4167000000002
$ws.Range("H96").Value2 = 1683.6666
$ws.Range("I96").Value2 = 1573.6666
$ws.Range("J96").Value2 = 1757
$ws.Range("K96").Value2 = 4720.9998
$ws.Range("L96").Value2 = 5271
$ws.Range("M96").Value2 = -3347.9998
$ws.Range("N96").Value2 = -8017
$ws.Range("H132").Value2 = 3713.484
$ws.Range("I132").Value2 = 3718.6072
$ws.Range("K132").Value2 = 11155.8216
$ws.Range("M132").Value2 = -8625.821599999999
$ws.Range("H137").Value2 = 4293.405
$ws.Range("I137").Value2 = 2474.0386
$ws.Range("J137").Value2 = 7249.875
$ws.Range("K137").Value2 = 7422.1158
$ws.Range("L137").Value2 = 21749.625
$ws.Range("M137").Value2 = -4872.1158
$ws.Range("N137").Value2 = -26849.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 3599.1
$ws.Range("J2").Value2 = 4270.2856
$ws.Range("L2").Value2 = 4270.2856
$ws.Range("N2").Value2 = -4496.2856
$ws.Range("H34").Value2 = 25137.5
$ws.Range("I34").Value2 = 20027
$ws.Range("K34").Value2 = 20027
$ws.Range("M34").Value2 = -19756
$ws.Range("H39").Value2 = 12000
$ws.Range("I39").Value2 = 15000
$ws.Range("J39").Value2 = 9000
$ws.Range("K39").Value2 = 15000
$ws.Range("L39").Value2 = 9000
$ws.Range("M39").Value2 = -14480
$ws.Range("N39").Value2 = -10040
$ws.Range("H74").Value2 = 77010940
$ws.Range("I74").Value2 = 166854700
$ws.Range("J74").Value2 = 2002.1428
$ws.Range("K74").Value2 = 166854700
$ws.Range("L74").Value2 = 2002.1428
$ws.Range("M74").Value2 = -166853826
$ws.Range("N74").Value2 = -3750.1428
$ws.Range("H77").Value2 = 77010940
$ws.Range("I77").Value2 = 166854700
$ws.Range("J77").Value2 = 2002.1428
$ws.Range("K77").Value2 = 834273500
$ws.Range("L77").Value2 = 10010.714
$ws.Range("M77").Value2 = -834269132
$ws.Range("N77").Value2 = -18746.714
$ws.Range("H116").Value2 = 3599.1
$ws.Range("J116").Value2 = 4270.2856
$ws.Range("L116").Value2 = 4270.2856
$ws.Range("N116").Value2 = -8858.285599999999
$ws.Range("H132").Value2 = 37140110
$ws.Range("I132").Value2 = 19366.889
$ws.Range("K132").Value2 = 58100.667
$ws.Range("M132").Value2 = -55570.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 3599.1
$ws.Range("J3").Value2 = 4270.2856
$ws.Range("L3").Value2 = 4270.2856
$ws.Range("N3").Value2 = -4498.2856
$ws.Range("H86").Value2 = 13018.311
$ws.Range("I86").Value2 = 6465.52
$ws.Range("J86").Value2 = 53973.25
$ws.Range("K86").Value2 = 6465.52
$ws.Range("L86").Value2 = 53973.25
$ws.Range("M86").Value2 = -5342.52
$ws.Range("N86").Value2 = -56219.25
$ws.Range("H89").Value2 = 13018.311
$ws.Range("I89").Value2 = 6465.52
$ws.Range("J89").Value2 = 53973.25
$ws.Range("K89").Value2 = 32327.6
$ws.Range("L89").Value2 = 269866.25
$ws.Range("M89").Value2 = -26711.6
$ws.Range("N89").Value2 = -281098.25
$ws.Range("H134").Value2 = 3253.3914
$ws.Range("I134").Value2 = 2616.65
$ws.Range("K134").Value2 = 7849.950000000001
$ws.Range("M134").Value2 = -5314.950000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value2 = 3035.2222
$ws.Range("J35").Value2 = 1200
$ws.Range("L35").Value2 = 1200
$ws.Range("N35").Value2 = -1788
$ws.Range("H86").Value2 = 7000.5
$ws.Range("I86").Value2 = 7000.5
$ws.Range("K86").Value2 = 7000.5
$ws.Range("M86").Value2 = -5877.5
$ws.Range("H89").Value2 = 7000.5
$ws.Range("I89").Value2 = 7000.5
$ws.Range("K89").Value2 = 35002.5
$ws.Range("M89").Value2 = -29386.5
$ws.Range("H94").Value2 = 1824.7931
$ws.Range("I94").Value2 = 1599.8
$ws.Range("J94").Value2 = 1871.6666
$ws.Range("K94").Value2 = 1599.8
$ws.Range("L94").Value2 = 1871.6666
$ws.Range("M94").Value2 = -1148.8
$ws.Range("N94").Value2 = -2773.6666
$ws.Range("H132").Value2 = 79734.766
$ws.Range("I132").Value2 = 79734.766
$ws.Range("K132").Value2 = 239204.298
$ws.Range("M132").Value2 = -236674.298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 1000
$ws.Range("I80").Value2 = 500
$ws.Range("J80").Value2 = 1500
$ws.Range("K80").Value2 = 1500
$ws.Range("L80").Value2 = 4500
$ws.Range("M80").Value2 = -564
$ws.Range("N80").Value2 = -6372
$ws.Range("H83").Value2 = 1000
$ws.Range("I83").Value2 = 500
$ws.Range("J83").Value2 = 1500
$ws.Range("K83").Value2 = 4500
$ws.Range("L83").Value2 = 13500
$ws.Range("M83").Value2 = 180
$ws.Range("N83").Value2 = -22860
$ws.Range("H87").Value2 = 899
$ws.Range("I87").Value2 = 899
$ws.Range("J87").Value2 = 0
$ws.Range("K87").Value2 = 2697
$ws.Range("L87").Value2 = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value2 = -1449
$ws.Range("H90").Value2 = 899
$ws.Range("I90").Value2 = 899
$ws.Range("J90").Value2 = 0
$ws.Range("K90").Value2 = 8091
$ws.Range("L90").Value2 = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value2 = -1851
$ws.Range("H141").Value2 = 4271
$ws.Range("I141").Value2 = 1240
$ws.Range("K141").Value2 = 3720
$ws.Range("M141").Value2 = 1460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value2 = 35000
$ws.Range("J18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("N18").ClearContents()
$ws.Range("H35").Value2 = 10000
$ws.Range("J35").Value2 = 10000
$ws.Range("L35").Value2 = 10000
$ws.Range("N35").Value2 = -10596
$ws.Range("H59").Value2 = 10000
$ws.Range("J59").Value2 = 10000
$ws.Range("L59").Value2 = 10000
$ws.Range("N59").Value2 = -11166
$ws.Range("H97").Value2 = 1897.3226
$ws.Range("I97").Value2 = 1451.1428
$ws.Range("J97").Value2 = 2834.3
$ws.Range("K97").Value2 = 1451.1428
$ws.Range("L97").Value2 = 2834.3
$ws.Range("M97").Value2 = -955.1428000000001
$ws.Range("N97").Value2 = -3826.3
$ws.Range("H113").Value2 = 3378.7354
$ws.Range("I113").Value2 = 2647.4443
$ws.Range("J113").Value2 = 4201.4375
$ws.Range("K113").Value2 = 2647.4443
$ws.Range("L113").Value2 = 4201.4375
$ws.Range("M113").Value2 = -477.4443000000001
$ws.Range("N113").Value2 = -8541.4375
$ws.Range("H122").Value2 = 2834.611
$ws.Range("I122").Value2 = 2787.6428
$ws.Range("K122").Value2 = 8362.928400000001
$ws.Range("M122").Value2 = -5912.928400000001
$ws.Range("H132").Value2 = 3180
$ws.Range("I132").Value2 = 3462.8
$ws.Range("K132").Value2 = 10388.4
$ws.Range("M132").Value2 = -7858.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1858.1
$ws.Range("I16").Value2 = 1618.5
$ws.Range("K16").Value2 = 1618.5
$ws.Range("M16").Value2 = -1448.5
$ws.Range("H22").Value2 = 2240.1667
$ws.Range("I22").Value2 = 1795.6666
$ws.Range("J22").Value2 = 2684.6667
$ws.Range("K22").Value2 = 1795.6666
$ws.Range("L22").Value2 = 2684.6667
$ws.Range("M22").Value2 = -1500.6666
$ws.Range("N22").Value2 = -3274.6667
$ws.Range("H27").Value2 = 2240.1667
$ws.Range("I27").Value2 = 1795.6666
$ws.Range("J27").Value2 = 2684.6667
$ws.Range("K27").Value2 = 1795.6666
$ws.Range("L27").Value2 = 2684.6667
$ws.Range("M27").Value2 = -1688.6666
$ws.Range("N27").Value2 = -2898.6667
$ws.Range("H32").Value2 = 0
$ws.Range("I32").Value2 = 0
$ws.Range("K32").Value2 = 0
$ws.Range("M32").ClearContents()
$ws.Range("H43").Value2 = 28509.334
$ws.Range("J43").Value2 = 29411.2
$ws.Range("L43").Value2 = 29411.2
$ws.Range("N43").Value2 = -29797.2
$ws.Range("H59").Value2 = 0
$ws.Range("J59").Value2 = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value2 = 0
$ws.Range("H82").Value2 = 2720.2354
$ws.Range("J82").Value2 = 2985.5
$ws.Range("L82").Value2 = 2985.5
$ws.Range("N82").Value2 = -3707.5
$ws.Range("H85").Value2 = 2720.2354
$ws.Range("J85").Value2 = 2985.5
$ws.Range("L85").Value2 = 2985.5
$ws.Range("N85").Value2 = -5481.5
$ws.Range("H98").Value2 = 80265
$ws.Range("J98").Value2 = 80265
$ws.Range("L98").Value2 = 80265
$ws.Range("N98").Value2 = -86255
$ws.Range("H99").Value2 = 38785.125
$ws.Range("I99").Value2 = 33199.4
$ws.Range("K99").Value2 = 33199.4
$ws.Range("M99").Value2 = -30204.4
$ws.Range("H122").Value2 = 4379.0713
$ws.Range("I122").Value2 = 3160.9
$ws.Range("J122").Value2 = 7424.5
$ws.Range("K122").Value2 = 9482.700000000001
$ws.Range("L122").Value2 = 22273.5
$ws.Range("M122").Value2 = -7032.700000000001
$ws.Range("N122").Value2 = -27173.5
$ws.Range("H132").Value2 = 76930810
$ws.Range("I132").Value2 = 4870.6523
$ws.Range("J132").Value2 = 666696300
$ws.Range("K132").Value2 = 14611.9569
$ws.Range("L132").Value2 = 2000088900
$ws.Range("M132").Value2 = -12081.9569
$ws.Range("N132").Value2 = -2000093960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value2 = 12971.5
$ws.Range("I40").Value2 = 7855.4
$ws.Range("K40").Value2 = 7855.4
$ws.Range("M40").Value2 = -7706.4
$ws.Range("J41").Value2 = 8239
$ws.Range("L41").Value2 = 8239
$ws.Range("N41").Value2 = -9019
$ws.Range("H107").Value2 = 2199.4736
$ws.Range("I107").Value2 = 1753.909
$ws.Range("J107").Value2 = 2812.125
$ws.Range("K107").Value2 = 5261.727000000001
$ws.Range("L107").Value2 = 8436.375
$ws.Range("M107").Value2 = -3341.727000000001
$ws.Range("N107").Value2 = -12276.375
$ws.Range("H122").Value2 = 38501948
$ws.Range("I122").Value2 = 45500850
$ws.Range("K122").Value2 = 136502550
$ws.Range("M122").Value2 = -136500100
